$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D18").Value = "2024년 2월 26일 동향: 모델 발전과 과학 혁신들"
$ws.Range("E18").Value = "https://freesearch.pe.kr/archives/5278"

$ws.Range("D28").Value = "[LLM] Transformer :: GPT 쓰긴 싫지만 GPT 안쓰면 도태될까봐 두려운 할미 MZ의 뒤늦은 LLM 공부 (1) Auto-regressive란?"
$ws.Range("E28").Value = "https://ropiens.tistory.com/263"

$wb.Save()
